$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 120 (pushes RAMP and everything below it down by one row)
# to add the new RALIGNXML command alphabetically after RALIGNC3D.
$ws.Rows.Item(120).Insert()

$ws.Cells.Item(120, 1).Value = "RALIGNXML"
$ws.Cells.Item(120, 2).Value = "Reads the alignment from a Land XML file"

# Update the active selection to reflect where the author's cursor ended up (B121)
$ws.Activate()
$ws.Range("B121").Select()
